$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (Source File) from 38 to 42 characters.
# Note: this runtime's ColumnWidth<->stored-width conversion has a constant
# +5/6 offset baked in (e.g. a stored width of 38 reads back as 37.17), so we
# compensate to land exactly on a stored width of 42.
$ws.Columns.Item(5).ColumnWidth = 41.166666666666664

# Update the source file name referenced in every data row (header is row 1)
$oldName = "pediatrics trial reference file.xlsx"
$newName = "pediatrics trial reference file (1).xlsx"

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 294) { $lastRow = 294 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq $oldName) {
        $cell.Value = $newName
    }
}
